$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)

for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
    $shp = $s.Shapes.Item($i)
    if ($shp.Name -eq "Picture 1") {
        $shp.Delete()
    }
}
